# This script applies the cryptocurrency price/volume update described by the
# commit "Updated cryptos list ... with GitHub Actions".
#
# - Most rows (2-39) only have their Price (column D) and Volume(1h) (column E)
#   text values refreshed to the latest scrape.
# - Row 40 (OKB) was dropped from the list, which shifted rows 41-51 up by one
#   position (Coin name, Link and Price updated accordingly), and a brand new
#   coin ("dogwifhat") was appended as the new row 51.
# - Row 42 (Bittensor) kept its position, only its Volume(1h) changed.
#
# Price/Volume cells in this sheet are stored as plain text (not numbers), even
# though several of the new values look like numbers (e.g. "531.51"). Setting
# such a string directly through .Value would make Excel auto-convert it to a
# real number, which would not match the original text-cell representation.
# SetText works around this by temporarily forcing the cell into text mode via
# a leading apostrophe, then restores the cell style so no stray numeric
# formatting is left behind.
function SetText($cell, $val) {
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.880.55'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '2.275.40'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.03%  '
SetText $ws.Range('D5') '531.51'
$ws.Range('E5').Value = '  -2.19%  '
SetText $ws.Range('D6') '130.59'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('E7').Value = '  +0.11%  '
SetText $ws.Range('D8') '0.583'
$ws.Range('E8').Value = '  +3.58%  '
$ws.Range('D9').Value = '2.274.27'
$ws.Range('E9').Value = '  +0.08%  '
SetText $ws.Range('D10') '0.0991'
$ws.Range('E10').Value = '  -1.71%  '
SetText $ws.Range('D11') '5.47'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('E12').Value = '  +0.80%  '
SetText $ws.Range('D13') '0.331'
$ws.Range('E13').Value = '  -0.71%  '
SetText $ws.Range('D14') '23.31'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').Value = '2.682.02'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').Value = '57.839.67'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '2.282.18'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  -1.37%  '
$ws.Range('E20').Value = '  -3.20%  '
SetText $ws.Range('D21') '311.89'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('E23').Value = '  +0.05%  '
SetText $ws.Range('D24') '62.32'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('E25').Value = '  -0.72%  '
SetText $ws.Range('D26') '0.997'
$ws.Range('E26').Value = '  -0.34%  '
SetText $ws.Range('D27') '7.95'
$ws.Range('E27').Value = '  -1.90%  '
SetText $ws.Range('D28') '1.26'
$ws.Range('E28').Value = '  -2.84%  '
SetText $ws.Range('D29') '170.15'
$ws.Range('E29').Value = '  -0.53%  '
SetText $ws.Range('D30') '1.70'
$ws.Range('E30').Value = '  -2.46%  '
$ws.Range('D31').Value = '0.0₃0717'
$ws.Range('E31').Value = '  +0.10%  '
SetText $ws.Range('D32') '5.74'
$ws.Range('E32').Value = '  -1.09%  '
SetText $ws.Range('D33') '1.05'
$ws.Range('E33').Value = '  -1.97%  '
SetText $ws.Range('D34') '0.377'
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  +0.00%  '
SetText $ws.Range('D36') '17.75'
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -2.04%  '
SetText $ws.Range('D39') '3.87'
$ws.Range('E39').Value = '  -1.87%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
SetText $ws.Range('D40') '1.48'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
SetText $ws.Range('D41') '139.54'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -2.85%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
SetText $ws.Range('D43') '3.41'
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
SetText $ws.Range('D44') '0.0945'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
SetText $ws.Range('D45') '0.0492'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
SetText $ws.Range('D46') '0.549'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
SetText $ws.Range('D47') '17.85'
$ws.Range('E47').Value = '  -2.18%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
SetText $ws.Range('D48') '0.0209'
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
SetText $ws.Range('D49') '10.95'
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('B50').Value = 'ZEEBU'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
SetText $ws.Range('D50') '4.64'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
SetText $ws.Range('D51') '1.51'
$ws.Range('E51').Value = '  +1.55%  '
